$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update DATE_TYPE_CODE (text field, keep leading zeros)
$ws.Range("J2").Value = "'004"
$ws.Range("J2").Style = "Normal"

# Update REPORT_DATE (text field storing a date-like string)
$ws.Range("N2").Value = "2020-09-30 00:00:00"

# Update numeric cash-flow figures
$ws.Range("O2").Value = 32188406.89
$ws.Range("P2").Value = 64.03901129979999
$ws.Range("Q2").Value = 453038220.54
$ws.Range("R2").Value = 901.3220139646
$ws.Range("S2").Value = 33412052.48
$ws.Range("T2").Value = 66.4734608839
$ws.Range("U2").Value = -45935252.32
$ws.Range("V2").Value = -91.3884353591
$ws.Range("Y2").Value = 45935252.32
$ws.Range("Z2").Value = 91.3884353591
$ws.Range("AA2").Value = -34276699.99
$ws.Range("AB2").Value = -68.1936818272
$ws.Range("AC2").Value = -50263747.42
$ws.Range("AD2").Value = 15.9828137696
